# Updated run for publication
# Updates the frequency-table values in Sheet1 (rows 2-5, columns B-X)
# to the refreshed values from the latest analysis run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.960869565217391
$ws.Range("C2").Value = 0.0268115942028985
$ws.Range("D2").Value = 0.794927536231884
$ws.Range("E2").Value = 0.426811594202899
$ws.Range("F2").Value = 0.938405797101449
$ws.Range("G2").Value = 0.00217391304347826
$ws.Range("H2").Value = 0.973913043478261
$ws.Range("I2").Value = 0.955797101449275
$ws.Range("J2").Value = 0.0108695652173913
$ws.Range("K2").Value = 0.0108695652173913
$ws.Range("L2").Value = 0.0920289855072464
$ws.Range("M2").Value = 0.00507246376811594
$ws.Range("N2").Value = 0.00434782608695652
$ws.Range("O2").Value = 0.951449275362319
$ws.Range("P2").Value = 0.201449275362319
$ws.Range("Q2").Value = 0.426811594202899
$ws.Range("R2").Value = 0.102898550724638
$ws.Range("S2").Value = 0.00434782608695652
$ws.Range("T2").Value = 0.022463768115942
$ws.Range("U2").Value = 0.00507246376811594
$ws.Range("V2").Value = 0.952173913043478
$ws.Range("W2").Value = 0.564492753623188
$ws.Range("X2").Value = 0.0528985507246377
$ws.Range("B3").Value = 0.0036231884057971
$ws.Range("C3").Value = 0.0181159420289855
$ws.Range("D3").Value = 0.0231884057971014
$ws.Range("E3").Value = 0.0123188405797101
$ws.Range("F3").Value = 0.0072463768115942
$ws.Range("G3").Value = 0.013768115942029
$ws.Range("H3").Value = 0.00942028985507246
$ws.Range("I3").Value = 0.0072463768115942
$ws.Range("J3").Value = 0.0036231884057971
$ws.Range("K3").Value = 0.0108695652173913
$ws.Range("L3").Value = 0.0289855072463768
$ws.Range("M3").Value = 0.144927536231884
$ws.Range("N3").Value = 0.0036231884057971
$ws.Range("P3").Value = 0.146376811594203
$ws.Range("Q3").Value = 0.0369565217391304
$ws.Range("R3").Value = 0.0297101449275362
$ws.Range("S3").Value = 0.00797101449275362
$ws.Range("T3").Value = 0.00289855072463768
$ws.Range("U3").Value = 0.00797101449275362
$ws.Range("V3").Value = 0.0173913043478261
$ws.Range("W3").Value = 0.0123188405797101
$ws.Range("X3").Value = 0.0159420289855072
$ws.Range("B4").Value = 0.0253623188405797
$ws.Range("C4").Value = 0.938405797101449
$ws.Range("D4").Value = 0.173188405797101
$ws.Range("E4").Value = 0.543478260869565
$ws.Range("F4").Value = 0.0376811594202899
$ws.Range("G4").Value = 0.981884057971015
$ws.Range("H4").Value = 0.00652173913043478
$ws.Range("I4").Value = 0.0268115942028985
$ws.Range("J4").Value = 0.96304347826087
$ws.Range("K4").Value = 0.0557971014492754
$ws.Range("L4").Value = 0.872463768115942
$ws.Range("M4").Value = 0.00869565217391304
$ws.Range("N4").Value = 0.0108695652173913
$ws.Range("O4").Value = 0.041304347826087
$ws.Range("P4").Value = 0.648550724637681
$ws.Range("Q4").Value = 0.00217391304347826
$ws.Range("R4").Value = 0.865217391304348
$ws.Range("S4").Value = 0.00144927536231884
$ws.Range("T4").Value = 0.969565217391304
$ws.Range("U4").Value = 0.986231884057971
$ws.Range("V4").Value = 0.013768115942029
$ws.Range("W4").Value = 0.406521739130435
$ws.Range("X4").Value = 0.913768115942029
$ws.Range("B5").Value = 0.0101449275362319
$ws.Range("C5").Value = 0.0159420289855072
$ws.Range("D5").Value = 0.00869565217391304
$ws.Range("E5").Value = 0.0159420289855072
$ws.Range("F5").Value = 0.0166666666666667
$ws.Range("G5").Value = 0.00217391304347826
$ws.Range("H5").Value = 0.0101449275362319
$ws.Range("I5").Value = 0.0101449275362319
$ws.Range("J5").Value = 0.022463768115942
$ws.Range("K5").Value = 0.922463768115942
$ws.Range("L5").Value = 0.00579710144927536
$ws.Range("M5").Value = 0.839855072463768
$ws.Range("N5").Value = 0.981159420289855
$ws.Range("O5").Value = 0.00652173913043478
$ws.Range("P5").Value = 0.00289855072463768
$ws.Range("Q5").Value = 0.534057971014493
$ws.Range("R5").Value = 0.00217391304347826
$ws.Range("S5").Value = 0.986231884057971
$ws.Range("T5").Value = 0.00434782608695652
$ws.Range("U5").Value = 0.00072463768115942
$ws.Range("V5").Value = 0.0166666666666667
$ws.Range("W5").Value = 0.0166666666666667
$ws.Range("X5").Value = 0.0166666666666667
